$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (rows 22-30), continuing the existing pattern.
$startRow = 22
$startA = 10002
$startB = 110021

for ($i = 0; $i -lt 9; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $startA + $i
    $ws.Range("B$r").Value = $startB + $i
    $ws.Range("C$r").Value = "eng"
    $ws.Range("D$r").Value = $true
    $ws.Range("E$r").Value = "superadmin"
    $ws.Range("F$r").Value = "now()"
}

# Select the rows below the data, matching the post-edit selection state.
$ws.Rows("31:1048576").Select() | Out-Null

# Set page setup (portrait orientation + print settings) for the sheet.
$ws.PageSetup.Orientation = 1
